$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.960.31"
$ws.Range("E2").Value = "  +3.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.862.30"
$ws.Range("E3").Value = "  +2.31%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.72"
$ws.Range("E5").Value = "  +2.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6395"
$ws.Range("E6").Value = "  +4.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9995"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3003"
$ws.Range("E8").Value = "  +4.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07494"
$ws.Range("E9").Value = "  +2.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.57"
$ws.Range("E10").Value = "  +7.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07682"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.870.39"
$ws.Range("E12").Value = "  +2.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.064"
$ws.Range("E13").Value = "  +2.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6929"
$ws.Range("E14").Value = "  +5.45%  "
$ws.Range("E15").Value = "  +3.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009426"
$ws.Range("E16").Value = "  +5.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.081"
$ws.Range("E17").Value = "  +4.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.927.65"
$ws.Range("E18").Value = "  +3.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.131.22"
$ws.Range("E19").Value = "  +3.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "239.68"
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("E21").Value = "  +2.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9997"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.397"
$ws.Range("E23").Value = "  +4.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.34"
$ws.Range("E25").Value = "  +1.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1425"
$ws.Range("E26").Value = "  +1.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.602"
$ws.Range("E27").Value = "  +2.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.01"
$ws.Range("E28").Value = "  +2.56%  "
$ws.Range("E29").Value = "  +1.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06066"
$ws.Range("E30").Value = "  +9.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.290"
$ws.Range("E31").Value = "  +7.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.140"
$ws.Range("E32").Value = "  +1.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.150"
$ws.Range("E33").Value = "  +1.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.889"
$ws.Range("E34").Value = "  +4.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.171"
$ws.Range("E35").Value = "  +3.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7309"
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.868"
$ws.Range("E38").Value = "  +1.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01809"
$ws.Range("E39").Value = "  +3.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.227.28"
$ws.Range("E40").Value = "  +2.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9303"
$ws.Range("E41").Value = "  +3.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.281"
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.036.79"
$ws.Range("E43").Value = "  +3.28%  "
$ws.Range("E44").Value = "  +0.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.37"
$ws.Range("E45").Value = "  +1.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.41"
$ws.Range("E46").Value = "  +3.17%  "
$ws.Range("E47").Value = "  +2.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5103"
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.358"
$ws.Range("E49").Value = "  +4.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4101"
$ws.Range("E50").Value = "  +2.88%  "
$ws.Range("E51").Value = "  +3.68%  "
